$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 6508
$ws1.Range("F4").Value = 750
$ws1.Range("F5").Value = 1094
$ws1.Range("F6").Value = 94
$ws1.Range("F7").Value = 576
$ws1.Range("F10").Value = 756
$ws1.Range("F19").Value = 689
$ws1.Range("F28").Value = 410
$ws1.Range("F32").Value = 668

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G18").Value = 680
$ws2.Range("F33").Value = 1683

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 1214
$ws3.Range("F7").Value = 438
$ws3.Range("F10").Value = 860

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1214
$ws4.Range("F6").Value = 438
$ws4.Range("F8").Value = 860
$ws4.Range("F9").Value = 6508
$ws4.Range("F11").Value = 750
$ws4.Range("F13").Value = 94
$ws4.Range("F14").Value = 576
$ws4.Range("F17").Value = 756
$ws4.Range("G28").Value = 680
$ws4.Range("F33").Value = 689
$ws4.Range("F43").Value = 1683
$ws4.Range("F44").Value = 1683
$ws4.Range("F46").Value = 410
$ws4.Range("F51").Value = 668
